# Apply the commit's data refresh to 上海-漫展信息.xlsx
# Sheets: 1=展览 (Exhibition), 2=演出 (Performance), 3=本地生活 (Local life), 4=全部类型 (All types)

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------
# Sheet "展览" (Exhibition) - F/G numeric refreshes only
# ---------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 1090
$ws1.Range("G2").Value = 30
$ws1.Range("F3").Value = 4700
$ws1.Range("F5").Value = 189
$ws1.Range("F6").Value = 1872
$ws1.Range("F7").Value = 49
$ws1.Range("F8").Value = 770
$ws1.Range("F12").Value = 1145
$ws1.Range("F14").Value = 828
$ws1.Range("F15").Value = 1858
$ws1.Range("F16").Value = 572
$ws1.Range("F17").Value = 530
$ws1.Range("F19").Value = 203
$ws1.Range("F20").Value = 17
$ws1.Range("F21").Value = 17
$ws1.Range("F23").Value = 1200
$ws1.Range("F24").Value = 613
$ws1.Range("F25").Value = 2526
$ws1.Range("F26").Value = 7
$ws1.Range("F27").Value = 297
$ws1.Range("F28").Value = 1602
$ws1.Range("F30").Value = 497
$ws1.Range("F33").Value = 4306

# ---------------------------------------------------------------
# Sheet "演出" (Performance) - F/G numeric refreshes
# ---------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F9").Value = 4166
$ws2.Range("F11").Value = 31
$ws2.Range("F17").Value = 286
$ws2.Range("G27").Value = 680
$ws2.Range("F28").Value = 1744

# ---------------------------------------------------------------
# Sheet "演出" (Performance) - insert a new event as row 30,
# pushing the existing rows 30-38 down to 31-39.
# ---------------------------------------------------------------

# Shift the existing data (columns B:I, rows 30-38) down by one row.
# Copy/paste (rather than Rows.Insert) keeps the original styles intact.
$ws2.Range("B30:I38").Copy()
$ws2.Range("B31:I39").PasteSpecial(-4104)
$ws2.Application.CutCopyMode = $false

# Give the new last row (39) the same "index" cell style as row 38 had.
$ws2.Range("A38").Copy()
$ws2.Range("A39").PasteSpecial(-4122)
$ws2.Application.CutCopyMode = $false
$ws2.Cells.Item(39, 1).Value = 38

# Update the dimension-driving bits: A30 keeps its original index (29),
# row 30 becomes the newly added event. B30 ("2024-05-03") is already the
# correct value (it is the same date as the event that used to sit in row 30
# and got pushed into row 31), so it is intentionally left untouched here to
# avoid Excel's autodetection turning the text date into a numeric date.
$ws2.Range("A30").Value = 29
$ws2.Range("C30").Value = "【大会员提前购】上海·申放送-Virtual Super Live-2024 in Shanghai"
$ws2.Range("D30").Value = "中兴路1599号金融街融泰中心 蜚声上海PHASE LIVE HOUSE"
$ws2.Range("E30").Value = "2024.05.03 19:00-05.03 22:00"
$ws2.Range("F30").Value = 7
$ws2.Range("G30").Value = 388
$ws2.Range("H30").Value = "https://show.bilibili.com/platform/detail.html?id=83102"
$ws2.Range("I30").Value = "//i1.hdslb.com/bfs/openplatform/202403/fkBE5inM1710929581281.png"

# ---------------------------------------------------------------
# Sheet "本地生活" (Local life) - F numeric refreshes
# ---------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("本地生活")
$ws3.Range("F5").Value = 1740
$ws3.Range("F6").Value = 1090
$ws3.Range("F7").Value = 318

# ---------------------------------------------------------------
# Sheet "全部类型" (All types) - F/G numeric refreshes
# ---------------------------------------------------------------
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F4").Value = 1740
$ws4.Range("F5").Value = 1090
$ws4.Range("F7").Value = 1090
$ws4.Range("G7").Value = 30
$ws4.Range("F9").Value = 4700
$ws4.Range("F11").Value = 189
$ws4.Range("F12").Value = 1872
$ws4.Range("F13").Value = 49
$ws4.Range("F14").Value = 770
$ws4.Range("F20").Value = 1145
$ws4.Range("F21").Value = 31
$ws4.Range("F25").Value = 828
$ws4.Range("F26").Value = 1858
$ws4.Range("F27").Value = 572
$ws4.Range("F28").Value = 530
$ws4.Range("F31").Value = 17
$ws4.Range("F33").Value = 286
$ws4.Range("F37").Value = 1200
$ws4.Range("F39").Value = 2526
$ws4.Range("F41").Value = 7
$ws4.Range("G43").Value = 680
$ws4.Range("F44").Value = 1744
$ws4.Range("F45").Value = 1602
$ws4.Range("F46").Value = 497
$ws4.Range("F49").Value = 4306
